$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: I1 = "I0", J1 = "IF" (same bold/border/center style as
# the rest of row 1, e.g. H1 "IP").
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New data cells: I2 = 8, J2 = 8 (plain numeric, default style like C2:H2).
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
